$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): two new, more recent report dates are being
# tracked ("Jun_15" / "Jun_17"). The two existing date headers slide two
# columns to the right (B->D, C->E) and the two newest dates take over
# columns B and C.
$oldB1 = $ws.Range("B1").Value2
$oldC1 = $ws.Range("C1").Value2

$ws.Range("E1").Value = $oldC1
$ws.Range("D1").Value = $oldB1
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"

# --- Data rows (2-27): column C (the detail for the oldest tracked date)
# moves two columns over to E, carrying its value and any highlight style
# along with it. The freshly vacated C and D columns are filled with the
# same "UN" placeholder already used in column B.
$ws.Range("C2:C27").Cut($ws.Range("E2:E27"))
$ws.Range("C2:D27").ClearFormats()
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Match the column widths already used for the data columns.
$ws.Columns("C").ColumnWidth = 7.1666667
$ws.Columns("D").ColumnWidth = 7.1666667
$ws.Columns("E").ColumnWidth = 7.1666667
